# [IMP] New data for test environment
# Updates the sample/demo data used by the sale_order.xlsx template:
#   - client_order_ref values (column D, rows 2/5/8) get refreshed
#     from 2020/2021-ish references to new 2022 references
#   - the two numeric "order date" codes in column D (rows 3/4) are
#     bumped from 21xxxx to 22xxxx
#   - the worksheet selection/cursor is reset back to A1
#   - the window tab ratio is nudged slightly (best effort)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- client_order_ref text values -----------------------------------
$ws.Range("D2").Value = "P1/2022/0001"
$ws.Range("D5").Value = "IT/22/004"
$ws.Range("D8").Value = "P1/2022/0007"

# --- numeric order-date style codes ----------------------------------
$ws.Range("D3").Value = 220123
$ws.Range("D4").Value = 22011214

# --- reset active selection from J6 back to A1 ------------------------
$ws.Range("A1").Select() | Out-Null

# --- cosmetic window setting (best effort; tab ratio 50% -> 60%) ------
$excel.ActiveWindow.TabRatio = 0.6
